$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "vel I gain" in J1
$ws.Range("J1").Value = "vel I gain"

# Row 2 (motor 1)
$ws.Range("H2").Value = 600
$ws.Range("I2").Value = 700
$ws.Range("J2").Value = 20

# Row 3 (motor 2)
$ws.Range("H3").Value = 850
$ws.Range("I3").Value = 700
$ws.Range("J3").Value = 14

# Row 4 (motor 3)
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 350
$ws.Range("J4").Value = 16

# Row 5 (motor 4)
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 350
$ws.Range("J5").Value = 16

# Row 6 (motor 5)
$ws.Range("H6").Value = 100
$ws.Range("J6").Value = 40

# Row 7 (motor 6)
$ws.Range("H7").Value = 100
$ws.Range("J7").Value = 40

# Row 8 (motor 7)
$ws.Range("J8").Value = 0

# Row 9 (motor 8)
$ws.Range("J9").Value = 0

# Row 10 (motor 9)
$ws.Range("J10").Value = 0

# Row 11 (motor 10)
$ws.Range("J11").Value = 0

# Row 12: leave a touched-but-empty cell in A12 (matches source round-trip)
$ws.Range("A12").Value = 0
$ws.Range("A12").ClearContents()

# Selection change
$ws.Range("C12").Select()
